# Releve 02 backend update: add the GTBank statement block (rows 8-9) and the
# extra "Titulaire" detail line in B3, per the new export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B3: extra detail line next to "Titulaire" -----------------------------
$ws.Range("B3").Value = "01-Mar-2017 O 01-Mar-2017 5,000.00 8,039.82 FBNMOBILE:MOHAMMED MUSA HASKE/LIBRACIN"

# --- Row 8: first statement line --------------------------------------------
$ws.Range("A8").Value = "21/91"
$ws.Range("B8").Value = "Account No: 012345678 - (66432/1/13/0) Currency :Naira Period: 01/Mar/2017 To 16/Mar/2017 Opening Balance: Trans Date Reference Value Date Debit Credit Balance Remarks TRANSFER BETWEEN CUSTOMERS 000016170301074834000150646887 01-Mar-2017 O 01-Mar-2017 5,000.00 8,039.82 FBNMOBILE:MOHAMMED MUSA HASKE/LIBRACIN 636239513140470886 OGUMGBU CHETACHI SIXTUS to MOHAMMED MUSA HASKE CASH WITHDRAWAL FROM OUR ATM -004273- - GTBank 1902 Ahmadu Bello Jos PLNG CASH WITHDRAWAL FROM OUR ATM -001885- - GTBank Jengre Road JOS 1 PLNG ATM/POS ACOUNT TO ACCOUNT TRANSFER -812087-- 000027102422;Acct Trsf:musaMohammedhaske NIBSS Instant Payment Outward 000013170306144649000009994924 USSD NIP Transfer from : 07034433706 TO ACCESS/TANSI ITAMAN REF:000013170306144649000009994924 COMMISSION 000013170306144649000009994924 USSD NIP Transfer from : 07034433706 TO 01-Mar-2017 3320004273 01-Mar-2017 5,000.00 03-Mar-2017 3310001885 03-Mar-2017 3,000.00 39.82 06-Mar-2017 9999812087 06-Mar-2017 20,000.00 20,039.82 06-Mar-2017. 0 06-Mar-2017 10,000.00 10,039.82 06-Mar-2017, 0. 06-Mar-2017 100.00 9,939.82 ,CCEcs/TANSI ITAMAN REF:000013170306144649000009994924 06-Mar-2017 i) 06-Mar-2017 5.00 9,934.82 VALUE ADDED TAX USSD NIP Transfer from : 0703 POS/WEB PURCHASE TRANSACTION -008430- - 08-Mar-2017 9999008430 07-Mar-2017 9,800.00 134,02/ ns Sea wien waratar tae FASTTRACK CASH DEPOSIT PINPAD-42220045-"

# C8 / C9 hold amounts that look numeric ("3,039.82", "7,000.00") but must be
# stored as plain text, matching the source export. Round-trip them through a
# temporary text number-format so Excel doesn't coerce them to numbers, then
# drop the format again so no residual cell style is left behind.
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "3,039.82"
$ws.Range("C8").ClearFormats()

$ws.Range("D8").Value = "Cr"

# --- Row 9: second statement line -------------------------------------------
$ws.Range("A9").Value = "32/87"
$ws.Range("B9").Value = "08-Mar-2017 0 08-Mar-2017 35,000.00 35,134.82 388 CASH DEPOSIT BY MOHAMMED MUSA HASKE 000636245691498822226338096987 NIBSS Instant Payment Outward 000013170309113623000010650577 TO FB/CHUKWUALUKA ILIEMENE REF:000013170309113623000010650577 COMMISSION 000013170309113623000010650577 09-Mar-2017 O 09-Mar-2017 28,134.82 09-Mar-2017 0 09-Mar-2017 100.00 28,034.82 TO FB/CHUKWUALUKA ILIEMENE REF:000013170309113623000010650577 09-Mar-2017 0 09-Mar-2017 5.00 28,029.82 VALUE ADDED TAX TO FB/CHUKWUALUKA ILIEMENE NIBSS Instant Payment Outward 000013170310063946000010750550 Smartzee (haske006@gmail.com) TO STERLING/ADESIYAN.COM LIMITED REF:000013170310063946000010750550 COMMISSION 000013170310063946000010750550 Smartzee (haske006@gmail.com) TO 10-Mar-2017 0 10-Mar-2017 26,770.00 1,259.82 eae saiincalinccael 100.00 1,159.82 STERLING/ADESIYAN.COM LIMITED REF:000013170310063946000010750550 10-Mar-2017 fo) 10-Mar-2017 5.00 1,154.82 VALUE ADDED TAX Smartzee (haske006@gmail.com) 154.2 POS/WEB PURCHASE TRANSACTION -023992- - 10-Mar-2017 9999023992 10-Mar-2017 1,000.00 ae ee ay ae wwrananaber Lane 13-Mar-2017 9999009251 «13-Mar-2017 10,000.00 10,154.82 FUNDS TRANSFER -009251- -FCMB Jos Br Jos PLNG TRANSFER BETWEEN CUSTOMERS via Internet 13-Mar-2017 fe) 13-Mar-2017 10,000.00 154.82 Banking from MOHAMMED MUSA HASKE to SUNDAY CHUKWUDI VICTOR 20,154.82 ATM/POS ACOUNT TO ACCOUNT TRANSFER -978298-- 13-Mar-2017 9999978298 13-Mar-2017 20,000.00 — 000027366089; Acct Trsf: NIBSS Instant Payment Outward 154.82 000013170312195807000011199244 Musa M haske 13-Mar-2017 0 13-Mar-2017 20,000.00 . TO DBN/NWANERI ANGUS NNAMDI REF:000013170312195807000011199244"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "7,000.00"
$ws.Range("C9").ClearFormats()

$ws.Range("D9").Value = "Dr"
